$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.209.50"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.908.83"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "0.7330"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("D6").Value = "243.81"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "0.3133"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "26.84"
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("D10").Value = "0.06912"
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("D11").Value = "0.7778"
$ws.Range("D12").Value = "0.07984"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "1.902.67"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "5.253"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "91.41"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").Value = "30.154.76"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "14.22"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "5.837"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "240.46"
$ws.Range("E19").Value = "  -6.96%  "
$ws.Range("D20").Value = "0.000007795"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "2.132.07"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "6.766"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "9.392"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").Value = "165.75"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "0.1270"
$ws.Range("E28").Value = "  -4.55%  "
$ws.Range("D29").Value = "2.087"
$ws.Range("E29").Value = "  -8.94%  "
$ws.Range("D30").Value = "1.547"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "4.079"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").Value = "0.05157"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "0.7427"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "2.755"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").Value = "0.01938"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "6.352"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("D41").Value = "74.71"
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("D42").Value = "0.4437"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "1.930"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "0.8355"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "101.05"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.596"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "9.790"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").Value = "37.50"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.043.99"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "945.36"
$ws.Range("E51").Value = "  -3.82%  "
